$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, expanded citation text replacing the old Emadian-only reference
# for the Amycolatopsis-related rows (369-380) that should keep a Ref value.
$newRef = "Teeraphatpornchai, T., Nakajima-Kambe, T., Shigeno-Akutsu, Y., Nakayama, M., Nomura, N., Nakahara, T., & Uchiyama, H. (2003). Isolation and characterization of a bacterium that degrades various polyester-based biodegradable plastics. Biotechnology letters, 25(1), 23-28., Sukkhum, S., Tokuyama, S., Tamura, T., & Kitpreechavanich, V. (2009). A novel poly (L-lactide) degrading actinomycetes isolated from Thai forest soil, phylogenic relationship and the enzyme characterization. The Journal of general and applied microbiology, 55(6), 459-467., Kim, M. N., & Park, S. T. (2010). Degradation of poly (L‐lactide) by a mesophilic bacterium. Journal of applied polymer science, 117(1), 67-74., A. Chomchoei, W. Pathom-Aree, A. Yokota, C. Kanongnuch, S. Lumyong Amycolatopsis thailandensis sp. nov., a poly(l-lactic acid)-degrading actinomycete, isolated from soil Int. J. Syst. Evol. Microbiol., 61 (2011), pp. 839-843, Penkhrue, W., Khanongnuch, C., Masaki, K., Pathom-aree, W., Punyodom, W., & Lumyong, S. (2015). Isolation and screening of biopolymer-degrading microorganisms from northern Thailand. [World Journal of Microbiology and Biotechnology, 31(9), 1431-1442.] via Emadian, S. M., Onay, T. T., & Demirel, B. (2017). Biodegradation of bioplastics in natural environments. Waste management, 59, 526-536."

for ($r = 369; $r -le 380; $r++) {
    $ws.Range("D$r").Value = $newRef
}

# Rows 381-489 lose their Ref (column D) value entirely.
for ($r = 381; $r -le 489; $r++) {
    $ws.Range("D$r").Value = ""
}

# Restore the cursor/selection position shown in the edited workbook.
$ws.Range("A381").Select()
